$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("F:F").Insert()
$ws.Range("F3").Value = "Is awarded"

$ws.Range("F3:F4").Merge()
$ws.Range("F3:F4").Select()

$shp = $ws.Shapes.Item(1)
$shp.Width = $ws.Cells.Item(1, 13).Left + 17
